$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value from 288 to 437
$ws.Range("B2").Value = 437

# Delete row 4 entirely (A4/B4 removed, dimension becomes A1:B3)
$ws.Rows.Item(4).Delete()
